$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 values ---
$ws.Range("G2").Value = 0.005389000000000001
$ws.Range("H2").Value = 0.016167
$ws.Range("M2").Value = 0.8377936666666667
$ws.Range("N2").Value = 2.513381
$ws.Range("O2").Value = 0.7130909380817101
$ws.Range("P2").Value = 0.7130909380817101
$ws.Range("Q2").Value = 0.004514870069666667
$ws.Range("R2").Value = 0.040633830627
$ws.Range("S2").Value = 0.7130909380817101
$ws.Range("T2").Value = 0.7130909380817101

# --- Add new row 3 (duplicate of row2's text columns, new numeric data) ---
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Pdyn"
$ws.Range("C3").Value = "Oprm1"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.005389000000000001
$ws.Range("H3").Value = 0.016167
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.3370826666666667
$ws.Range("N3").Value = 1.011248
$ws.Range("O3").Value = 0.2869090619182899
$ws.Range("P3").Value = 0.2869090619182899
$ws.Range("Q3").Value = 0.001816538490666667
$ws.Range("R3").Value = 0.016348846416
$ws.Range("S3").Value = 0.2869090619182899
$ws.Range("T3").Value = 0.2869090619182899
